$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Formatting first: give the brand-new cells (column G entirely, plus the
#    lower half of column F) the same style as their already-formatted
#    neighbours, BEFORE any values are written, so the later .Value
#    assignments don't have to fight with freshly-created default styles.
# ---------------------------------------------------------------------------
$ws.Range("E2:E9").Copy()
$ws.Range("G2:G9").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("E5:E9").Copy()
$ws.Range("F5:F9").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)      # xlPasteFormats (header look)
$excel.CutCopyMode = 0

# Column F needs to be a bit wider to fit the new "gastou hoje" entries.
$ws.Range("F1").ColumnWidth = 17.166666666666668   # -> stored width 18

# ---------------------------------------------------------------------------
# 2) Values - Planejamento (existing rows, tweaked) + Execução (new PESSOA(S),
#    "gastou hoje" and "Total" columns) for the PERT/Gantt/burndown tables.
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = "1h + 1h + 2h"
$ws.Range("F3").Value = "50min + 2h"
$ws.Range("E4").Value = "Ana e Beatriz / Marco"
$ws.Range("E6").Value = "Marco e Laís"
$ws.Range("E8").Value = "Vinícius, Laís e Jader"
$ws.Range("E7").Value = "Próxima Sprint"
$ws.Range("F8").Value = "1h "
$ws.Range("F7").Value = "0h"
$ws.Range("G7").Value = "0h"
$ws.Range("F4").Value = "50min + 2h + 1h10"
$ws.Range("E9").Value = "Ana, Beatriz e Laís"
$ws.Range("F6").Value = "2h + 3h + 1h"
$ws.Range("G9").Value = "3h"
$ws.Range("G1").Value = "Total"
$ws.Range("G2").Value = " 4h"
$ws.Range("G3").Value = " 3h"
$ws.Range("G4").Value = "4h"
$ws.Range("G6").Value = "6h"
$ws.Range("F9").Value = "1h + 2h"
$ws.Range("E5").Value = "Jader"
$ws.Range("F5").Value = "1h"
$ws.Range("G5").Value = "1h"
$ws.Range("G8").Value = "1h"

# ---------------------------------------------------------------------------
# 3) Selection follows the last cell the author was working on.
# ---------------------------------------------------------------------------
$ws.Range("F9").Select()
